$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# Remove the "Disclaimer" sheet entirely.
$null = $wb.Worksheets("Disclaimer").Delete()

# The "White_List" sheet loses its "QDM Category" column (old column E),
# which shifts every subsequent column left by one.
$ws = $wb.Worksheets("White_List")
$ws.Activate()
$null = $ws.Range("E1").EntireColumn.Delete()

# Restore the selection recorded on the White_List sheet.
$null = $ws.Range("G32").Select()
